# Update "想去人数" (F column) counts for both the "展览" sheet and the
# "全部类型" sheet, which both contain the same set of events (the latter
# being an aggregated view that also includes rows from "演出").
$wb = $excel.ActiveWorkbook

$updates1 = @{
    2  = 476
    4  = 8010
    9  = 115
    10 = 469
    15 = 78
    17 = 5906
    18 = 187
    19 = 274
    20 = 1907
    21 = 15
    22 = 27
    24 = 406
}

$updates4 = @{
    2  = 476
    4  = 8010
    9  = 115
    10 = 469
    15 = 78
    18 = 5906
    20 = 187
    21 = 274
    22 = 1907
    23 = 15
    24 = 27
    26 = 406
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates1.Keys) {
    $ws1.Range("F$row").Value = $updates1[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
